$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = $ws.Cells.Item($row, 2).Style
}

Set-TextValue $ws 2 4 "42.925.68"
Set-TextValue $ws 2 5 "  -5.42%  "
Set-TextValue $ws 3 4 "2.210.84"
Set-TextValue $ws 3 5 "  -6.70%  "
Set-TextValue $ws 4 5 "  +0.33%  "
Set-TextValue $ws 5 4 "314.75"
Set-TextValue $ws 5 5 "  +1.34%  "
Set-TextValue $ws 6 4 "97.74"
Set-TextValue $ws 6 5 "  -9.97%  "
Set-TextValue $ws 7 4 "0.579"
Set-TextValue $ws 7 5 "  -7.82%  "
Set-TextValue $ws 8 5 "  +0.22%  "
Set-TextValue $ws 9 4 "0.557"
Set-TextValue $ws 9 5 "  -9.68%  "
Set-TextValue $ws 10 4 "36.47"
Set-TextValue $ws 10 5 "  -11.50%  "
Set-TextValue $ws 11 4 "54.30"
Set-TextValue $ws 11 5 "  -2.00%  "
Set-TextValue $ws 12 4 "0.0827"
Set-TextValue $ws 12 5 "  -10.15%  "
Set-TextValue $ws 13 4 "7.74"
Set-TextValue $ws 13 5 "  -9.05%  "
Set-TextValue $ws 14 4 "0.105"
Set-TextValue $ws 14 5 "  -4.07%  "
Set-TextValue $ws 15 4 "0.862"
Set-TextValue $ws 15 5 "  -12.22%  "
Set-TextValue $ws 16 4 "2.551.30"
Set-TextValue $ws 16 5 "  -6.35%  "
Set-TextValue $ws 17 4 "14.05"
Set-TextValue $ws 17 5 "  -7.90%  "
Set-TextValue $ws 18 4 "2.211.43"
Set-TextValue $ws 18 5 "  -6.15%  "
Set-TextValue $ws 19 4 "42.798.79"
Set-TextValue $ws 19 5 "  -5.47%  "
Set-TextValue $ws 20 4 "14.57"
Set-TextValue $ws 20 5 "  +1.42%  "
Set-TextValue $ws 21 4 "0.0₃0954"
Set-TextValue $ws 21 5 "  -10.03%  "
Set-TextValue $ws 22 4 "6.39"
Set-TextValue $ws 22 5 "  -12.77%  "
Set-TextValue $ws 23 4 "65.13"
Set-TextValue $ws 23 5 "  -11.06%  "
Set-TextValue $ws 24 4 "3.18"
Set-TextValue $ws 24 5 "  -9.24%  "
Set-TextValue $ws 25 4 "236.72"
Set-TextValue $ws 25 5 "  -9.19%  "
Set-TextValue $ws 26 4 "2.12"
Set-TextValue $ws 26 5 "  -8.63%  "
Set-TextValue $ws 27 5 "  -0.16%  "
Set-TextValue $ws 28 4 "10.05"
Set-TextValue $ws 28 5 "  -10.15%  "
Set-TextValue $ws 29 4 "2.23"
Set-TextValue $ws 29 5 "  -5.03%  "
Set-TextValue $ws 30 4 "6.24"
Set-TextValue $ws 30 5 "  -15.17%  "
Set-TextValue $ws 31 4 "0.0880"
Set-TextValue $ws 31 5 "  -9.12%  "
Set-TextValue $ws 32 4 "20.40"
Set-TextValue $ws 32 5 "  -8.81%  "
Set-TextValue $ws 33 4 "33.68"
Set-TextValue $ws 33 5 "  -10.64%  "
Set-TextValue $ws 34 4 "154.68"
Set-TextValue $ws 34 5 "  -8.56%  "
Set-TextValue $ws 35 4 "2.77"
Set-TextValue $ws 35 5 "  -5.63%  "
Set-TextValue $ws 36 4 "3.21"
Set-TextValue $ws 36 5 "  +8.12%  "
Set-TextValue $ws 37 4 "1.99"
Set-TextValue $ws 37 5 "  +13.82%  "
Set-TextValue $ws 38 5 "  -6.62%  "
Set-TextValue $ws 39 4 "4.41"
Set-TextValue $ws 39 5 "  -7.89%  "
Set-TextValue $ws 40 4 "0.103"
Set-TextValue $ws 40 5 "  -12.22%  "
Set-TextValue $ws 41 4 "3.68"
Set-TextValue $ws 41 5 "  -6.11%  "
Set-TextValue $ws 42 4 "0.0323"
Set-TextValue $ws 42 5 "  -9.04%  "
Set-TextValue $ws 43 4 "1.865.00"
Set-TextValue $ws 43 5 "  +11.58%  "
Set-TextValue $ws 44 5 "  +0.38%  "
Set-TextValue $ws 45 4 "12.27"
Set-TextValue $ws 45 5 "  -5.18%  "
Set-TextValue $ws 46 4 "88.90"
Set-TextValue $ws 46 5 "  -10.67%  "
Set-TextValue $ws 47 4 "0.206"
Set-TextValue $ws 47 5 "  -11.52%  "
Set-TextValue $ws 48 4 "5.42"
Set-TextValue $ws 48 5 "  -2.41%  "
Set-TextValue $ws 49 4 "76.27"
Set-TextValue $ws 49 5 "  -5.90%  "
Set-TextValue $ws 50 4 "59.91"
Set-TextValue $ws 50 5 "  -14.08%  "
Set-TextValue $ws 51 4 "8.62"
Set-TextValue $ws 51 5 "  -6.77%  "
